$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "invest"
$ws.Range("C2").Value = -0.8214

# Row 4
$ws.Range("B4").Value = "uncertain"
$ws.Range("C4").Value = -0.2036

# Row 7
$ws.Range("B7").Value = "inflation"
$ws.Range("C7").Value = -0.4892

# Row 8
$ws.Range("B8").Value = "trade"
$ws.Range("C8").Value = 0.1004

# Row 9
$ws.Range("B9").Value = "interest"
$ws.Range("C9").Value = -0.4244

# Row 10
$ws.Range("B10").Value = "uncertain"
$ws.Range("C10").Value = -0.3641

# Row 11
$ws.Range("B11").Value = "invest"
$ws.Range("C11").Value = 0.0982

# Row 12
$ws.Range("B12").Value = "trade"
$ws.Range("C12").Value = 0.8056

# Row 13
$ws.Range("B13").Value = "uncertain"
$ws.Range("C13").Value = -0.3662

# Row 14
$ws.Range("B14").Value = "interest"
$ws.Range("C14").Value = 0.4306

# Row 15
$ws.Range("B15").Value = "invest"
$ws.Range("C15").Value = 0.4222

# Row 17
$ws.Range("B17").Value = "trade"
$ws.Range("C17").Value = -0.6305

# Row 18
$ws.Range("B18").Value = "inflation"
$ws.Range("C18").Value = 0.0717

# Row 19
$ws.Range("B19").Value = "interest"
$ws.Range("C19").Value = 0.0838

# Row 20
$ws.Range("B20").Value = "invest"
$ws.Range("C20").Value = -0.0479

# Row 21
$ws.Range("B21").Value = "uncertain"
$ws.Range("C21").Value = -0.557
